$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows whose ID is "RM 232" and "SC 92" (entire rows, shifting the rest up)
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# After the deletions, the remaining rows (SC 5, SC 101, SC 105, SC 119, SC 120, SC 132, SC 193, SC 232)
# occupy rows 26-33. Update column D ("C" header) values to match the new masking pattern.
$ws.Range("D27").Value = -14.6      # SC 101
$ws.Range("D28").ClearContents()    # SC 105 -> blank
$ws.Range("D29").ClearContents()    # SC 119 -> blank
$ws.Range("D30").Value = -13.6      # SC 120
$ws.Range("D32").ClearContents()    # SC 193 -> blank
